$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ONSITE AM")
Write-Host $ws.Name
